$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# --- Update the reporting period / update dates stored in row 8 ---
# Year of the report: 2021 -> 2022
$ws.Range("A8").Value = 2022
# Period start: 2021-10-01 -> 2022-01-01
$ws.Range("B8").Value = 44562
# Period end: 2021-12-31 -> 2022-03-31
$ws.Range("C8").Value = 44651
# "Fecha de actualización" / "Fecha de validación": 2022-01-10 -> 2022-04-08
$ws.Range("AJ8").Value = 44659
$ws.Range("AK8").Value = 44659

# --- Update the view: scroll back to the left of the sheet and select C8 ---
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("C8").Select()

# --- Narrower workbook window width ---
$win.Width = 15600
